$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213 (shifts existing rows 213-297 down to 214-298)
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new record
$ws.Cells.Item(213, 1).Value = 1
$ws.Cells.Item(213, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(213, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(213, 4).Value = 44825
$ws.Cells.Item(213, 5).Value = 15
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100108
$ws.Cells.Item(213, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(213, 9).Value = 100108006
$ws.Cells.Item(213, 10).Value = "Plátano"
$ws.Cells.Item(213, 11).Value = "Sin especificar"
$ws.Cells.Item(213, 12).Value = "Pintón"
$ws.Cells.Item(213, 13).Value = 120
$ws.Cells.Item(213, 14).Value = 24000
$ws.Cells.Item(213, 15).Value = 25000
$ws.Cells.Item(213, 16).Value = 24500
$ws.Cells.Item(213, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(213, 18).Value = "Ecuador"
$ws.Cells.Item(213, 19).Value = 1225
$ws.Cells.Item(213, 20).Value = 20
